$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Valor_Observado" column (old column E); this shifts
# old F:N left into E:M and updates the used-range dimension automatically.
$ws.Range("E:E").Delete()

# Rename / shorten header labels (B1:D1) and relabel the model columns (E1:M1).
$ws.Range('B1').Value = 'Config'
$ws.Range('C1').Value = 'Dist'
$ws.Range('D1').Value = 'Var'
$ws.Range('E1').Value = 'Block Bootstrapping'
$ws.Range('F1').Value = 'Sieve Bootstrap'
$ws.Range('G1').Value = 'LSPM'
$ws.Range('H1').Value = 'LSPMW'
$ws.Range('I1').Value = 'AREPD'
$ws.Range('J1').Value = 'MCPS'
$ws.Range('K1').Value = 'AV-MCPS'
$ws.Range('L1').Value = 'DeepAR'
$ws.Range('M1').Value = 'EnCQR-LSTM'

# Replace the model-output values (E2:M25) with the new simulation results.
$ws.Range('E2').Value = 0.8302605154986817
$ws.Range('F2').Value = 0.5701484598298238
$ws.Range('G2').Value = 0.6154830006432697
$ws.Range('H2').Value = 0.6010188200491328
$ws.Range('I2').Value = 0.58632834667427
$ws.Range('J2').Value = 1.209135245155832
$ws.Range('K2').Value = 1.684266228803616
$ws.Range('L2').Value = 0.5808960493311638
$ws.Range('M2').Value = 1.398368984411266
$ws.Range('E3').Value = 0.9421339880699815
$ws.Range('F3').Value = 0.6008207339982492
$ws.Range('G3').Value = 1.720662475744687
$ws.Range('H3').Value = 1.4140443721926
$ws.Range('I3').Value = 1.641200821511794
$ws.Range('J3').Value = 0.7090579793933162
$ws.Range('K3').Value = 0.6049685694849769
$ws.Range('L3').Value = 0.5838194819475898
$ws.Range('M3').Value = 1.26239589110202
$ws.Range('E4').Value = 0.7641699819969016
$ws.Range('F4').Value = 0.5684726953474861
$ws.Range('G4').Value = 0.7434725158980061
$ws.Range('H4').Value = 1.033924083173391
$ws.Range('I4').Value = 1.212516135575344
$ws.Range('J4').Value = 0.5891487775665455
$ws.Range('K4').Value = 0.5951847356001096
$ws.Range('L4').Value = 0.5668684924877472
$ws.Range('M4').Value = 1.261790439022001
$ws.Range('E5').Value = 0.7001349100041439
$ws.Range('F5').Value = 0.5801201865876834
$ws.Range('G5').Value = 0.5592357420085849
$ws.Range('H5').Value = 0.9581844708557178
$ws.Range('I5').Value = 1.120519337928296
$ws.Range('J5').Value = 0.631872421268074
$ws.Range('K5').Value = 0.5902739862936782
$ws.Range('L5').Value = 0.5666291809639427
$ws.Range('M5').Value = 1.265052063648472
$ws.Range('E6').Value = 0.694294991394674
$ws.Range('F6').Value = 0.5716273156115166
$ws.Range('G6').Value = 0.581667387843658
$ws.Range('H6').Value = 0.735580343151773
$ws.Range('I6').Value = 0.8497178670746374
$ws.Range('J6').Value = 0.6122053229718336
$ws.Range('K6').Value = 0.6203770068090583
$ws.Range('L6').Value = 0.5958858431234384
$ws.Range('M6').Value = 1.293413932768203
$ws.Range('E7').Value = 0.9980337098506654
$ws.Range('F7').Value = 0.5862149665744185
$ws.Range('G7').Value = 1.294980686999779
$ws.Range('H7').Value = 0.7484228292778472
$ws.Range('I7').Value = 0.6757624538003401
$ws.Range('J7').Value = 0.797913236001338
$ws.Range('K7').Value = 1.435862002060615
$ws.Range('L7').Value = 0.6007579253455961
$ws.Range('M7').Value = 1.449970708376357
$ws.Range('E8').Value = 0.7374500403794718
$ws.Range('F8').Value = 0.6040112872230934
$ws.Range('G8').Value = 1.69306634651521
$ws.Range('H8').Value = 0.9535067215038631
$ws.Range('I8').Value = 1.111212433967447
$ws.Range('J8').Value = 0.6230663043770373
$ws.Range('K8').Value = 0.5810702541662629
$ws.Range('L8').Value = 0.5763278065159443
$ws.Range('M8').Value = 1.260844169126452
$ws.Range('E9').Value = 0.8564915296420436
$ws.Range('F9').Value = 0.5556391435194928
$ws.Range('G9').Value = 0.5627252797330612
$ws.Range('H9').Value = 1.31365232046561
$ws.Range('I9').Value = 1.51529452768268
$ws.Range('J9').Value = 0.5735342617516785
$ws.Range('K9').Value = 0.6432641734995187
$ws.Range('L9').Value = 0.5614781726673669
$ws.Range('M9').Value = 1.252658236408253
$ws.Range('E10').Value = 0.7065114556136874
$ws.Range('F10').Value = 0.6321098104107928
$ws.Range('G10').Value = 0.7182021712251789
$ws.Range('H10').Value = 0.8746253959144359
$ws.Range('I10').Value = 1.011034051274231
$ws.Range('J10').Value = 0.6743437732980865
$ws.Range('K10').Value = 1.53619128290221
$ws.Range('L10').Value = 0.6053339363823407
$ws.Range('M10').Value = 1.271147362012382
$ws.Range('E11').Value = 1.337945280234004
$ws.Range('F11').Value = 0.5878837798846194
$ws.Range('G11').Value = 1.323585182868722
$ws.Range('H11').Value = 1.908525830852693
$ws.Range('I11').Value = 2.14088406004404
$ws.Range('J11').Value = 1.282550014139896
$ws.Range('K11').Value = 0.7868334672641082
$ws.Range('L11').Value = 0.5769383694053177
$ws.Range('M11').Value = 1.295094959575668
$ws.Range('E12').Value = 0.6704673294315043
$ws.Range('F12').Value = 0.5837012135520754
$ws.Range('G12').Value = 1.299660959898055
$ws.Range('H12').Value = 0.779131913559705
$ws.Range('I12').Value = 0.8899388950224162
$ws.Range('J12').Value = 0.6356407054224767
$ws.Range('K12').Value = 0.632780246899402
$ws.Range('L12').Value = 0.621222357670189
$ws.Range('M12').Value = 1.281080105357782
$ws.Range('E13').Value = 0.6555951148093688
$ws.Range('F13').Value = 0.5689322459107138
$ws.Range('G13').Value = 0.5567080543239267
$ws.Range('H13').Value = 0.6285371712945618
$ws.Range('I13').Value = 0.7043249013662766
$ws.Range('J13').Value = 0.6247212384416525
$ws.Range('K13').Value = 0.5892047237395046
$ws.Range('L13').Value = 0.584455538069996
$ws.Range('M13').Value = 1.305218844637772
$ws.Range('E14').Value = 0.5680212950708862
$ws.Range('F14').Value = 0.5652833450043651
$ws.Range('G14').Value = 0.5710401231320216
$ws.Range('H14').Value = 0.6061563679200483
$ws.Range('I14').Value = 0.6001245969432121
$ws.Range('J14').Value = 0.7455175628346088
$ws.Range('K14').Value = 0.5919898165403417
$ws.Range('L14').Value = 0.5651574663645219
$ws.Range('M14').Value = 0.8320518246777392
$ws.Range('E15').Value = 0.875992783500889
$ws.Range('F15').Value = 0.5925744158154862
$ws.Range('G15').Value = 0.8877029160656006
$ws.Range('H15').Value = 1.103866121509306
$ws.Range('I15').Value = 1.050261054027521
$ws.Range('J15').Value = 0.8100745295542449
$ws.Range('K15').Value = 1.002627431354512
$ws.Range('L15').Value = 0.5800838093383122
$ws.Range('M15').Value = 0.9296437207440696
$ws.Range('E16').Value = 0.5583036942403639
$ws.Range('F16').Value = 0.5574337545168069
$ws.Range('G16').Value = 0.705937904167238
$ws.Range('H16').Value = 0.599161686507844
$ws.Range('I16').Value = 0.5959471877453024
$ws.Range('J16').Value = 0.5980348690087427
$ws.Range('K16').Value = 0.6024202091150952
$ws.Range('L16').Value = 0.5674232886348218
$ws.Range('M16').Value = 0.8234115892635129
$ws.Range('E17').Value = 0.7698895614609779
$ws.Range('F17').Value = 0.5834856991616761
$ws.Range('G17').Value = 0.6837411957009446
$ws.Range('H17').Value = 0.9097673115665995
$ws.Range('I17').Value = 0.8596341298864951
$ws.Range('J17').Value = 0.6846995001093747
$ws.Range('K17').Value = 0.6119649904108392
$ws.Range('L17').Value = 0.574456705107848
$ws.Range('M17').Value = 0.8666892049492608
$ws.Range('E18').Value = 0.5844641969886624
$ws.Range('F18').Value = 0.557692952208632
$ws.Range('G18').Value = 0.7862293453891511
$ws.Range('H18').Value = 0.5647466833689541
$ws.Range('I18').Value = 0.593148336758103
$ws.Range('J18').Value = 0.5951106155961495
$ws.Range('K18').Value = 0.6815542396200416
$ws.Range('L18').Value = 0.5682614174776449
$ws.Range('M18').Value = 0.8532866110038613
$ws.Range('E19').Value = 1.135862819602279
$ws.Range('F19').Value = 0.5873440740797472
$ws.Range('G19').Value = 1.003914114579771
$ws.Range('H19').Value = 1.053156359433132
$ws.Range('I19').Value = 1.234396452718643
$ws.Range('J19').Value = 0.7669532996888468
$ws.Range('K19').Value = 0.7663790766495802
$ws.Range('L19').Value = 0.6607703233438025
$ws.Range('M19').Value = 1.161504424599149
$ws.Range('E20').Value = 0.5786470439461949
$ws.Range('F20').Value = 0.5679888250900044
$ws.Range('G20').Value = 1.354605221417615
$ws.Range('H20').Value = 0.6495423808652246
$ws.Range('I20').Value = 0.624585388231011
$ws.Range('J20').Value = 0.8864275608603933
$ws.Range('K20').Value = 0.7340269466151524
$ws.Range('L20').Value = 0.6481023637948258
$ws.Range('M20').Value = 0.8292558635767765
$ws.Range('E21').Value = 0.6579713632432502
$ws.Range('F21').Value = 0.5746238287168732
$ws.Range('G21').Value = 0.6975945571235654
$ws.Range('H21').Value = 0.6135005959019243
$ws.Range('I21').Value = 0.6778500635979854
$ws.Range('J21').Value = 0.6691347036245856
$ws.Range('K21').Value = 0.7157151239804969
$ws.Range('L21').Value = 0.5846201446191676
$ws.Range('M21').Value = 0.8909507439593609
$ws.Range('E22').Value = 0.6399601709646265
$ws.Range('F22').Value = 0.5644225862264284
$ws.Range('G22').Value = 0.5587800646371598
$ws.Range('H22').Value = 0.6052500672640517
$ws.Range('I22').Value = 0.6697287701594986
$ws.Range('J22').Value = 0.6052037984860295
$ws.Range('K22').Value = 0.6481457380999851
$ws.Range('L22').Value = 0.5665041426884475
$ws.Range('M22').Value = 0.8842317881207914
$ws.Range('E23').Value = 0.7833636431385682
$ws.Range('F23').Value = 0.5724078507470147
$ws.Range('G23').Value = 1.037358730651962
$ws.Range('H23').Value = 0.9593846783702812
$ws.Range('I23').Value = 0.9072159212282638
$ws.Range('J23').Value = 1.239697626299048
$ws.Range('K23').Value = 1.124686572161028
$ws.Range('L23').Value = 0.594937671181363
$ws.Range('M23').Value = 0.8837676183690185
$ws.Range('E24').Value = 0.6239388765475302
$ws.Range('F24').Value = 0.5825362015112325
$ws.Range('G24').Value = 0.9315112342100603
$ws.Range('H24').Value = 0.5999528404946788
$ws.Range('I24').Value = 0.6468242894697617
$ws.Range('J24').Value = 0.7062808102267715
$ws.Range('K24').Value = 0.5949423705926189
$ws.Range('L24').Value = 0.5820539449133243
$ws.Range('M24').Value = 0.8656543066564143
$ws.Range('E25').Value = 0.6026291714254723
$ws.Range('F25').Value = 0.5783931498863421
$ws.Range('G25').Value = 0.771960384807206
$ws.Range('H25').Value = 0.6955395109800564
$ws.Range('I25').Value = 0.668991725446905
$ws.Range('J25').Value = 0.9755333613876803
$ws.Range('K25').Value = 0.6892277635285358
$ws.Range('L25').Value = 0.5953036887671866
$ws.Range('M25').Value = 0.8377748525000902
